$wb = $excel.ActiveWorkbook

# --- Sprint2Info: Utilization input changed from 0.14 to 0.13 ---
$infoWs = $wb.Worksheets.Item("Sprint2Info")
$infoWs.Range("B8").Value = 0.13

# --- BurnDown2Table: Actual Burn Down (column D) manual values updated ---
$burnWs = $wb.Worksheets.Item("BurnDown2Table")
$burnWs.Range("D4").Value = 91
$burnWs.Range("D5").Value = 87
$burnWs.Range("D6").Value = 80
$burnWs.Range("D7").Value = 76
$burnWs.Range("D8").Value = 73
$burnWs.Range("D9").Value = 70
$burnWs.Range("D10").Value = 66
$burnWs.Range("D11").Value = 61
$burnWs.Range("D13").Value = 55
$burnWs.Range("D14").Value = 50
$burnWs.Range("D15").Value = 45
$burnWs.Range("D16").Value = 40
$burnWs.Range("D17").Value = 37
$burnWs.Range("D18").Value = 33
$burnWs.Range("D19").Value = 27
$burnWs.Range("D20").Value = 25
$burnWs.Range("D21").Value = 22
$burnWs.Range("D22").Value = 15
$burnWs.Range("D23").Value = 10
$burnWs.Range("D24").Value = 5
$burnWs.Range("D25").Value = 0

$burnWs.Range("D26").Select() | Out-Null
$burnWs.Activate()
